$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume data per the latest scrape.
# Rows 18/19, 24/25, 29/30, 39/40 also swap ranking order (coin name,
# link, price, volume all change), and row 51 replaces FLOKI with Monero.
# A leading apostrophe forces numeric-looking price strings to stay text,
# matching the source sheet where every Price/Volume cell is stored as text.

$ws.Range("D2").Value = '71.783.42'
$ws.Range("E2").Value = '  +4.52%  '
$ws.Range("D3").Value = '4.023.81'
$ws.Range("E3").Value = '  +4.53%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = '''536.25'
$ws.Range("E5").Value = '  +3.51%  '
$ws.Range("D6").Value = '''153.52'
$ws.Range("E6").Value = '  +9.10%  '
$ws.Range("D7").Value = '''0.693'
$ws.Range("E7").Value = '  +14.21%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '''0.753'
$ws.Range("E9").Value = '  +5.90%  '
$ws.Range("D10").Value = '''0.173'
$ws.Range("E10").Value = '  +2.89%  '
$ws.Range("D11").Value = '''0.0000328'
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").Value = '''48.29'
$ws.Range("E12").Value = '  +16.07%  '
$ws.Range("D13").Value = '''10.80'
$ws.Range("E13").Value = '  +4.85%  '
$ws.Range("D14").Value = '4.665.42'
$ws.Range("E14").Value = '  +4.77%  '
$ws.Range("D15").Value = '4.031.51'
$ws.Range("E15").Value = '  +3.92%  '
$ws.Range("D16").Value = '''14.21'
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("D17").Value = '''20.57'
$ws.Range("E17").Value = '  -2.71%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '''0.133'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("B19").Value = 'Polygon'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D19").Value = '''1.20'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("D20").Value = '71.683.53'
$ws.Range("E20").Value = '  +4.60%  '
$ws.Range("D21").Value = '''432.56'
$ws.Range("E21").Value = '  +4.38%  '
$ws.Range("D22").Value = '''99.08'
$ws.Range("E22").Value = '  +14.41%  '
$ws.Range("D23").Value = '''3.54'
$ws.Range("E23").Value = '  +2.41%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '''4.22'
$ws.Range("E24").Value = '  +5.60%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").Value = '''14.55'
$ws.Range("E25").Value = '  +3.88%  '
$ws.Range("D26").Value = '''11.12'
$ws.Range("E26").Value = '  -9.08%  '
$ws.Range("D27").Value = '''10.85'
$ws.Range("E27").Value = '  +4.42%  '
$ws.Range("D28").Value = '''3.72'
$ws.Range("E28").Value = '  +30.70%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").Value = '''5.84'
$ws.Range("E29").Value = '  +2.92%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '''36.95'
$ws.Range("E30").Value = '  +4.58%  '
$ws.Range("D31").Value = '''13.49'
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("D32").Value = '''0.131'
$ws.Range("E32").Value = '  +5.42%  '
$ws.Range("D33").Value = '''685.72'
$ws.Range("E33").Value = '  +1.26%  '
$ws.Range("D34").Value = '''6.94'
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("D35").Value = '''66.76'
$ws.Range("E35").Value = '  +0.38%  '
$ws.Range("D36").Value = '''42.36'
$ws.Range("E36").Value = '  +7.45%  '
$ws.Range("D37").Value = '''0.427'
$ws.Range("E37").Value = '  -4.43%  '
$ws.Range("D38").Value = '''0.156'
$ws.Range("E38").Value = '  +6.03%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0830'
$ws.Range("E39").Value = '  -2.12%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '''3.48'
$ws.Range("E40").Value = '  +11.09%  '
$ws.Range("D41").Value = '''3.45'
$ws.Range("E41").Value = '  +2.29%  '
$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").Value = '''0.999'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '''0.0489'
$ws.Range("E44").Value = '  +3.29%  '
$ws.Range("D45").Value = '''0.151'
$ws.Range("E45").Value = '  +6.74%  '
$ws.Range("D46").Value = '''2.66'
$ws.Range("E46").Value = '  -7.59%  '
$ws.Range("D47").Value = '''3.40'
$ws.Range("E47").Value = '  -0.94%  '
$ws.Range("D48").Value = '''9.55'
$ws.Range("E48").Value = '  +9.28%  '
$ws.Range("D49").Value = '''3.03'
$ws.Range("E49").Value = '  +0.86%  '
$ws.Range("D50").Value = '''3.34'
$ws.Range("E50").Value = '  +1.54%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D51").Value = '''143.44'
$ws.Range("E51").Value = '  +0.41%  '
